# Applies the "adds custom units, troubleshoots make_metadata" edit to the
# feather_metadata.xlsx workbook.
#
# Summary of the change:
#  - dataset!A2        -> dataset name filled in
#  - personnel!A2:E2   -> a new personnel row (Ryon Kurth, project lead) filled in,
#                         which also retires the stray "creator" placeholder text
#                         that used to sit in personnel!D2
#  - keyword_set       -> new taxon keywords appended (chinook already present;
#                         Speckled dace, Steelhead trout [x2], Tule perch added)
#  - taxonomic_coverage-> same new taxa appended as new rows
#  - coverage!F2:G2    -> begin_date / end_date filled in
#  - various sheets    -> selection / active-cell bookkeeping

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Cell content edits (order matters: it controls the order new strings are
#    appended to the shared-string table, so we add the brand-new values in
#    the same sequence the diff shows them appearing).
# ---------------------------------------------------------------------------

$wsKeywords = $wb.Worksheets.Item("keyword_set")
$wsTaxo     = $wb.Worksheets.Item("taxonomic_coverage")
$wsPersonnel = $wb.Worksheets.Item("personnel")
$wsDataset   = $wb.Worksheets.Item("dataset")
$wsCoverage  = $wb.Worksheets.Item("coverage")

# keyword_set: existing "chinook" row, then the new taxa
$wsKeywords.Range("A8").Value = "chinook"
$wsKeywords.Range("A9").Value = "Speckled dace"
$wsKeywords.Range("A10").Value = "Steelhead trout "
$wsKeywords.Range("A11").Value = "Steelhead trout"
$wsKeywords.Range("A12").Value = "Tule perch"

# taxonomic_coverage: same new taxa, as new rows (chinook already in A2)
$wsTaxo.Range("A3").Value = "Speckled dace"
$wsTaxo.Range("A4").Value = "Steelhead trout "
$wsTaxo.Range("A5").Value = "Steelhead trout"
$wsTaxo.Range("A6").Value = "Tule perch"

# personnel: fill in the second personnel row (replaces the stray "creator"
# text that used to live in D2)
$wsPersonnel.Range("A2").Value = "Ryon"
$wsPersonnel.Range("B2").Value = "Kurth"
$wsPersonnel.Range("C2").Value = "Ryon.Kurth@water.ca.gov "
$wsPersonnel.Range("D2").Value = "project lead"

# dataset: fill in the dataset name
$wsDataset.Range("A2").Value = " Feather River mini-Snorkel surveying"

# coverage: fill in begin_date / end_date (stored as date serials, existing
# date-format style on these cells is left untouched)
$wsCoverage.Range("F2").Value = 36963
$wsCoverage.Range("G2").Value = 37125

# ---------------------------------------------------------------------------
# 2. Selection / active-cell bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------

$wb.Worksheets.Item("personnel").Range("D3").Select()
$wb.Worksheets.Item("title").Range("C8").Select()
$wb.Worksheets.Item("keyword_set").Range("A8:A12").Select()
$wb.Worksheets.Item("license").Range("A4").Select()
$wb.Worksheets.Item("taxonomic_coverage").Range("A7:A11").Select()

# "coverage" becomes the active sheet/tab, with E2 selected - do this last so
# it ends up the active (tabSelected) sheet on save.
$wsCoverage.Range("E2").Select()
$wsCoverage.Activate()
